$d = $word.ActiveDocument

# --- 1) First paragraph: pad the existing sentence with two trailing
#        spaces, then append a new (separate) run in dark red holding
#        the "(This is a change ... )" annotation. ---
$p1 = $d.Paragraphs(1)
$firstRange = $p1.Range
$textOnly = $d.Range($firstRange.Start, $firstRange.End - 1)
$textOnly.InsertAfter("  ")

$p1 = $d.Paragraphs(1)
$firstRange = $p1.Range
$newRun = $d.Range($firstRange.End - 1, $firstRange.End - 1)
$newRun.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch alternate)")
$newRun.Font.Color = 192

# --- 2) Insert one extra empty paragraph right before the very last
#        (already empty) paragraph of the document, i.e. between the
#        "Normal (Web)" empty paragraph and the final empty paragraph
#        that precedes the section properties. ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
